$p = $ppt.ActivePresentation

# 1) Slide 3 ("Nuestra propuesta de trabajo"): append " y sus respuestas" to the
#    last bullet's trailing run text.
$s3 = $p.Slides.Item(3)
$shp = $s3.Shapes.Item(2)
$tf = $shp.TextFrame.TextRange
$para = $tf.Paragraphs($tf.Paragraphs().Count, 1)
$runCount = $para.Runs().Count
$lastRun = $para.Runs($runCount, 1)
$lastRun.Text = ", también en formato HTML, con ejercicios y sus respuestas"

# 2) Remove slide 5 ("Sobre el uso de la inteligencia artificial") entirely.
$p.Slides.Item(5).Delete()
